# Apply the edits described by the diff:
#  - Fix typo "targt" -> "target" in cell A3 on both the "arousal" and
#    "valence" sheets.
#  - Append an asterisk to specific "< .05" p-value cells to mark them as
#    significant after switching from Bonferroni to FDR correction.

$wb = $excel.ActiveWorkbook

# ---- Sheet: arousal ----
$ws1 = $wb.Worksheets.Item("arousal")
$ws1.Range("A3").Value = "target"

$sheet1Cells = @("E12","H12","E14","H14","E16","H16","H20","H24","E26","H26","H28","H36","H40","H60","H70","E82","H82","H84","H88","E96","H96","E98","H98","E100","H100","E108","H108","E110","E112","H112","E120","E122","E128","H128")
foreach ($addr in $sheet1Cells) {
    $ws1.Range($addr).Value = "< .05*"
}

# ---- Sheet: valence ----
$ws2 = $wb.Worksheets.Item("valence")
$ws2.Range("A3").Value = "target"

$sheet2Cells = @("H5","H7","H9")
foreach ($addr in $sheet2Cells) {
    $ws2.Range($addr).Value = "< .05*"
}
